$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new data rows (alphabetically-sorted new pharmacy items) ---
# Current data rows occupy 7-16; totals row is 17; footer row is 18.
# We need 5 more data rows inserted right before the totals row (17),
# pushing totals -> 22 and footer -> 23.
$ws.Rows("17:21").Insert(-4121, 0)

# Copy the formatting (styles/borders/merges) of an existing data row (16)
# into each of the 5 newly-inserted blank rows so they look identical to
# the rest of the table.
for ($r = 17; $r -le 21; $r++) {
    $ws.Range("A16:Q16").Copy($ws.Range("A" + $r + ":Q" + $r))
}

# --- Write out the full (now 15-row) item table in alphabetical order ---
$rows = @(
    @{ A=1;  C="  1+1 INFINITY CARE 50+ SPF LOTION";    H="0:0";    L="0"; N="264.00"; P="264.0000"; Q="1:0" },
    @{ A=2;  C="ALKAPRESS PLUS 10/160MG 20 F.C. TABS."; H="0:1";    L="1"; N="102.00"; P="102.0000"; Q="1:0" },
    @{ A=3;  C="ANTI-COX II 15MG/3ML 6 AMP";             H="0:1";    L="1"; N="78.00";  P="12.4800";  Q="0:1" },
    @{ A=4;  C="DIAMICRON 60MG M.R. 30 SCORED TAB";      H="3:1";    L="1"; N="108.00"; P="108.0000"; Q="1:0" },
    @{ A=5;  C="DOLIPRANE 1 GM 15 TABS.";                 H="4:1";    L="1"; N="48.00";  P="15.8400";  Q="0:1" },
    @{ A=6;  C="EMPACOZA TRIO XR 25/5/1000  30TAB";       H="1:1";    L="0"; N="396.00"; P="130.6800"; Q="0:1" },
    @{ A=7;  C="ERASTAPEX 20 MG 30 F.C.TAB.";             H="1:0";    L="1"; N="75.00";  P="24.7500";  Q="0:1" },
    @{ A=8;  C="HIBIOTIC N 600MG SUSP. 80 ML";            H="1:0";    L="1"; N="92.00";  P="92.0000";  Q="1:0" },
    @{ A=9;  C="METACARDIA MR 35MG 30 F.C. TAB.";         H="0:2";    L="1"; N="60.00";  P="60.0000";  Q="1:0" },
    @{ A=10; C="MUCOSOL PED. 125MG/5ML SYRUP 120ML";      H="0:0";    L="1"; N="23.00";  P="23.0000";  Q="1:0" },
    @{ A=11; C="TERRAMYCIN EYE OINT. 5 GM";               H="2:0";    L="1"; N="28.00";  P="28.0000";  Q="1:0" },
    @{ A=12; C="WATER FOR INJECTION AMP. 5 ML";           H="8522:0"; L="1"; N="2.00";   P="4.0000";   Q="2:0" },
    @{ A=13; C="ZYROVAZET 10/20MG 30 F.C. TABLETS";       H="0:2";    L="1"; N="294.00"; P="97.0200";  Q="0:1" },
    @{ A=14; C="ايزي سويت قطعه";                           H="4:0";    L="0"; N="10.00";  P="20.0000";  Q="2:0" },
    @{ A=15; C="سرنجات 3 سم";                              H="0:0";    L="0"; N="2.00";   P="2.0000";   Q="1:0" }
)

$r = 7
foreach ($row in $rows) {
    $ws.Range("A" + $r).Value = $row.A
    $ws.Range("C" + $r).Value = $row.C
    $ws.Range("H" + $r).Value = $row.H
    $ws.Range("L" + $r).Value = $row.L
    $ws.Range("N" + $r).Value = $row.N
    $ws.Range("P" + $r).Value = $row.P
    $ws.Range("Q" + $r).Value = $row.Q
    $r = $r + 1
}

# --- Update the grand total (sum of the "sale price" column) ---
$ws.Range("P22").Value = 983.76999999999998

# --- Update the generated-on timestamp shown in the footer ---
$ws.Range("A23").Value = "Monday, 21 July, 2025 10:12 AM"
